# ctrl-settings\BEPEfCT updating to match 3.4.4 format
#
# This script reproduces, via Excel COM automation, the edits that were
# made to "Boolean Exempt Process Emissions from Carbon Tax.xlsx":
#   1. On the "About" sheet: rows 9/10 swap back to their correct order,
#      and two new explanatory rows (13/14) are appended.
#   2. On the "BEPEfCT" sheet: the old single "Boolean" setting row is
#      replaced by a per-industry-sector breakdown (25 sector rows),
#      a new italic "Unit: boolean (0 or 1)" label is added in A1, the
#      first column is widened to fit the longer labels, and the page
#      is explicitly set to portrait orientation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Rows 9 and 10 swap places (their text content is exchanged).
$oldA9  = $about.Cells.Item(9, 1).Text
$oldA10 = $about.Cells.Item(10, 1).Text
$about.Cells.Item(9, 1).Value  = $oldA10
$about.Cells.Item(10, 1).Value = $oldA9

# Two new rows are appended, explaining the U.S. carve-outs.
$about.Cells.Item(13, 1).Value = "In the U.S., we exempt agriculture and water and waste process emissions. Generally, "
$about.Cells.Item(14, 1).Value = "proposed taxes do not cover these sectors."

# ---------------------------------------------------------------------
# 2) "BEPEfCT" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BEPEfCT")

# A1 becomes a new italic unit label (previously held "Boolean").
$ws.Cells.Item(1, 1).Value = "Unit: boolean (0 or 1)"
$ws.Cells.Item(1, 1).Font.Italic = $true
# B1 keeps its existing title text ("Exempt Process Emissions from Carbon Tax").

$sectorNames = @(
    "agriculture and forestry 01T03",
    "coal mining 05",
    "oil and gas extraction 06",
    "other mining and quarrying 07T08",
    "food beverage and tobacco 10T12",
    "textiles apparel and leather 13T15",
    "wood products 16",
    "pulp paper and printing 17T18",
    "refined petroleum and coke 19",
    "chemicals 20",
    "rubber and plastic products 22",
    "glass and glass products 231",
    "cement and other nonmetallic minerals 239",
    "iron and steel 241",
    "other metals 242",
    "metal products except machinery and vehicles 25",
    "computers and electronics 26",
    "appliances and electrical equipment 27",
    "other machinery 28",
    "road vehicles 29",
    "nonroad vehicles 30",
    "other manufacturing 31T33",
    "energy pipelines and gas processing 352T353",
    "water and waste 36T39",
    "construction 41T43"
)

# Boolean flag per sector: only agriculture and water/waste are exempted.
$sectorValues = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $sectorNames.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $sectorNames[$i]
    $ws.Cells.Item($row, 2).Value = $sectorValues[$i]
}

# Widen column A so the longer sector labels are readable (matches the
# workbook's new custom width for column A, ~47.18 characters).
$ws.Columns.Item(1).ColumnWidth = 46.25

# Explicitly set portrait page orientation (adds a <pageSetup> element).
$ws.PageSetup.Orientation = 1
